$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save a snapshot of the current (pre-edit) D,K,L,M,N,O,P,Q,R,S,T values
# for each affected row, then redistribute them according to the cycles
# observed in the diff:
#   (2 -> 12 -> 8 -> 2)
#   (3 -> 15 -> 13 -> 7 -> 4 -> 3)
#   (9 -> 14 -> 9)
$rows = @(2, 3, 4, 7, 8, 9, 12, 13, 14, 15)
$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Mapping: destination row -> source row (i.e. destination receives source's data)
$destFromSource = @{
    2  = 8
    12 = 2
    8  = 12
    3  = 4
    15 = 3
    13 = 15
    7  = 13
    4  = 7
    9  = 14
    14 = 9
}

foreach ($dest in $destFromSource.Keys) {
    $src = $destFromSource[$dest]
    $srcData = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$dest").Value2 = $srcData[$c]
    }
}
